$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value2 = 6
$ws.Cells.Item(2,6).Value2 = 7.5
$ws.Cells.Item(2,7).Value2 = 0.306218628951098
$ws.Cells.Item(2,8).Value2 = 1.859452291935418
$ws.Cells.Item(2,9).Value2 = 0.7500793906648423
$ws.Cells.Item(2,10).Value2 = 0.2434902368758151
$lCell = $ws.Cells.Item(2,12)
$lCell.Value2 = "Small"
$lCell.Interior.Color = 13228184
$lCell.Font.Color = 16777215
$ws.Cells.Item(2,13).Value2 = 14.1
$nCell = $ws.Cells.Item(2,14)
$nCell.Value2 = 2.833333333333333
$nCell.Interior.Color = 6194490
$nCell.Font.Color = 15856113

# Row 3
$ws.Cells.Item(3,5).Value2 = 6
$ws.Cells.Item(3,6).Value2 = 4.666666666666667
$ws.Cells.Item(3,7).Value2 = 0.3630534566317305
$ws.Cells.Item(3,8).Value2 = 1.163940968374505
$ws.Cells.Item(3,9).Value2 = 0.8892957181014012
$ws.Cells.Item(3,10).Value2 = 0.2072883292816573
$lCell = $ws.Cells.Item(3,12)
$lCell.Value2 = "Small"
$lCell.Interior.Color = 13228184
$lCell.Font.Color = 16777215
$ws.Cells.Item(3,13).Value2 = 7.9
$nCell = $ws.Cells.Item(3,14)
$nCell.Value2 = 3.166666666666667
$nCell.Interior.Color = 15790831
$nCell.Font.Color = 0

# Row 4
$ws.Cells.Item(4,5).Value2 = 6
$ws.Cells.Item(4,6).Value2 = 2
$ws.Cells.Item(4,7).Value2 = -0.4649795257917584
$ws.Cells.Item(4,8).Value2 = 1.157353053076757
$ws.Cells.Item(4,9).Value2 = -1.138962579031098
$ws.Cells.Item(4,10).Value2 = 0.153164179321003
$lCell = $ws.Cells.Item(4,12)
$lCell.Value2 = "Small"
$lCell.Interior.Color = 13228184
$lCell.Font.Color = 16777215
$ws.Cells.Item(4,13).Value2 = 7.3
$nCell = $ws.Cells.Item(4,14)
$nCell.Value2 = 3.166666666666667
$nCell.Interior.Color = 15790831
$nCell.Font.Color = 0

# Row 5
$ws.Cells.Item(5,5).Value2 = 6
$ws.Cells.Item(5,6).Value2 = 4.5
$ws.Cells.Item(5,7).Value2 = 0.5083882833194932
$ws.Cells.Item(5,8).Value2 = 1.162272426319431
$ws.Cells.Item(5,9).Value2 = 1.245291885342247
$ws.Cells.Item(5,10).Value2 = 0.134100706077421
$lCell = $ws.Cells.Item(5,12)
$lCell.Value2 = "Medium"
$lCell.Interior.Color = 10732133
$lCell.Font.Color = 16777215
$ws.Cells.Item(5,13).Value2 = 6.4
$nCell = $ws.Cells.Item(5,14)
$nCell.Value2 = 3.166666666666667
$nCell.Interior.Color = 15790831
$nCell.Font.Color = 0

# Row 6
$ws.Cells.Item(6,5).Value2 = 6
$ws.Cells.Item(6,6).Value2 = 0
$ws.Cells.Item(6,7).Value2 = -1.240215946204668
$ws.Cells.Item(6,8).Value2 = 0.9052317076000181
$ws.Cells.Item(6,9).Value2 = -3.037896239064466
$ws.Cells.Item(6,10).Value2 = 0.01440999036614293
$lCell = $ws.Cells.Item(6,12)
$lCell.Value2 = "Very large"
$lCell.Interior.Color = 4491810
$lCell.Font.Color = 16777215
$ws.Cells.Item(6,13).Value2 = 4.9
$nCell = $ws.Cells.Item(6,14)
$nCell.Value2 = 3.5
$nCell.Interior.Color = 4602842
$nCell.Font.Color = 15856113

# Row 7
$ws.Cells.Item(7,5).Value2 = 6
$ws.Cells.Item(7,6).Value2 = 1.333333333333333
$ws.Cells.Item(7,7).Value2 = -0.4206017080997648
$ws.Cells.Item(7,8).Value2 = 1.353772275459665
$ws.Cells.Item(7,9).Value2 = -1.030259569787458
$ws.Cells.Item(7,10).Value2 = 0.1750615400789786
$lCell = $ws.Cells.Item(7,12)
$lCell.Value2 = "Small"
$lCell.Interior.Color = 13228184
$lCell.Font.Color = 16777215
$ws.Cells.Item(7,13).Value2 = 7.9
$nCell = $ws.Cells.Item(7,14)
$nCell.Value2 = 3.166666666666667
$nCell.Interior.Color = 15790831
$nCell.Font.Color = 0

# Row 8
$ws.Cells.Item(8,5).Value2 = 6
$ws.Cells.Item(8,6).Value2 = 3.833333333333333
$ws.Cells.Item(8,7).Value2 = 0.3301002535153649
$ws.Cells.Item(8,8).Value2 = 1.377772101762393
$ws.Cells.Item(8,9).Value2 = 0.8085771850760132
$ws.Cells.Item(8,10).Value2 = 0.2277467504197668
$lCell = $ws.Cells.Item(8,12)
$lCell.Value2 = "Small"
$lCell.Interior.Color = 13228184
$lCell.Font.Color = 16777215
$ws.Cells.Item(8,13).Value2 = 7.8
$nCell = $ws.Cells.Item(8,14)
$nCell.Value2 = 3.333333333333333
$nCell.Interior.Color = 10196454
$nCell.Font.Color = 15856113

# Row 9
$ws.Cells.Item(9,5).Value2 = 6
$ws.Cells.Item(9,6).Value2 = 3.166666666666667
$ws.Cells.Item(9,7).Value2 = 0.1915713475706197
$ws.Cells.Item(9,8).Value2 = 1.04683201543518
$ws.Cells.Item(9,9).Value2 = 0.469252050885384
$ws.Cells.Item(9,10).Value2 = 0.3293196862112076
$lCell = $ws.Cells.Item(9,12)
$lCell.Value2 = "Very small"
$lCell.Interior.Color = 15133900
$lCell.Font.Color = 16777215
$ws.Cells.Item(9,13).Value2 = 6
$nCell = $ws.Cells.Item(9,14)
$nCell.Value2 = 3.333333333333333
$nCell.Interior.Color = 10196454
$nCell.Font.Color = 15856113

# Row 10
$ws.Cells.Item(10,5).Value2 = 6
$ws.Cells.Item(10,6).Value2 = 0.8333333333333334
$ws.Cells.Item(10,7).Value2 = -0.6222704840087872
$ws.Cells.Item(10,8).Value2 = 1.110702196658986
$ws.Cells.Item(10,9).Value2 = -1.524245167816248
$ws.Cells.Item(10,10).Value2 = 0.09397796480394355
$lCell = $ws.Cells.Item(10,12)
$lCell.Value2 = "Medium"
$lCell.Interior.Color = 10732133
$lCell.Font.Color = 16777215
$ws.Cells.Item(10,13).Value2 = 5.8
$nCell = $ws.Cells.Item(10,14)
$nCell.Value2 = 3.333333333333333
$nCell.Interior.Color = 10196454
$nCell.Font.Color = 15856113

# Row 11
$ws.Cells.Item(11,5).Value2 = 6
$ws.Cells.Item(11,6).Value2 = 4.166666666666667
$ws.Cells.Item(11,7).Value2 = 0.6158551018724515
$ws.Cells.Item(11,8).Value2 = 1.157161714610398
$ws.Cells.Item(11,9).Value2 = 1.508530755077259
$ws.Cells.Item(11,10).Value2 = 0.095895221785726
$lCell = $ws.Cells.Item(11,12)
$lCell.Value2 = "Medium"
$lCell.Interior.Color = 10732133
$lCell.Font.Color = 16777215
$ws.Cells.Item(11,13).Value2 = 7.4
$nCell = $ws.Cells.Item(11,14)
$nCell.Value2 = 3.166666666666667
$nCell.Interior.Color = 15790831
$nCell.Font.Color = 0

# Row 12
$ws.Cells.Item(12,5).Value2 = 6
$ws.Cells.Item(12,6).Value2 = 3.5
$ws.Cells.Item(12,7).Value2 = 0.413788994367455
$ws.Cells.Item(12,8).Value2 = 1.141648520577947
$ws.Cells.Item(12,9).Value2 = 1.013571897379647
$ws.Cells.Item(12,10).Value2 = 0.1786474632144996
$lCell = $ws.Cells.Item(12,12)
$lCell.Value2 = "Small"
$lCell.Interior.Color = 13228184
$lCell.Font.Color = 16777215
$ws.Cells.Item(12,13).Value2 = 7.6
$nCell = $ws.Cells.Item(12,14)
$nCell.Value2 = 3.333333333333333
$nCell.Interior.Color = 10196454
$nCell.Font.Color = 15856113

# Row 13
$ws.Cells.Item(13,5).Value2 = 6
$ws.Cells.Item(13,6).Value2 = 5.5
$ws.Cells.Item(13,7).Value2 = 1.093538012305202
$ws.Cells.Item(13,8).Value2 = 1.178930328231524
$ws.Cells.Item(13,9).Value2 = 2.678610144485098
$ws.Cells.Item(13,10).Value2 = 0.0219465755021444
$lCell = $ws.Cells.Item(13,12)
$lCell.Value2 = "Large"
$lCell.Interior.Color = 7712064
$lCell.Font.Color = 16777215
$ws.Cells.Item(13,13).Value2 = 6.9
$nCell = $ws.Cells.Item(13,14)
$nCell.Value2 = 2.833333333333333
$nCell.Interior.Color = 6194490
$nCell.Font.Color = 15856113

# Row 14
$ws.Cells.Item(14,5).Value2 = 6
$ws.Cells.Item(14,6).Value2 = 1.666666666666667
$ws.Cells.Item(14,7).Value2 = -0.2433396078043851
$ws.Cells.Item(14,8).Value2 = 0.9636728036088147
$ws.Cells.Item(14,9).Value2 = -0.5960578733297227
$ws.Cells.Item(14,10).Value2 = 0.2885464810659479
$lCell = $ws.Cells.Item(14,12)
$lCell.Value2 = "Small"
$lCell.Interior.Color = 13228184
$lCell.Font.Color = 16777215
$ws.Cells.Item(14,13).Value2 = 5
$nCell = $ws.Cells.Item(14,14)
$nCell.Value2 = 3.333333333333333
$nCell.Interior.Color = 10196454
$nCell.Font.Color = 15856113
